$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row number, then D/E/G new values ($null = leave unchanged).
# All four data columns (D,E,G) are plain-text cells in this sheet (prices,
# percentages and the 'Hora' counter are stored as text, not numbers), so
# every write below is prefixed with a literal apostrophe to force Excel to
# keep it as text instead of auto-converting to a number/percentage/date.
$rows = @(
    @{ Row = 2; D = '257.79'; E = '0.04%'; G = '2' }
    @{ Row = 3; D = '26.94'; E = '-1.49%'; G = '2' }
    @{ Row = 4; D = '4.585'; E = '-12.09%'; G = '2' }
    @{ Row = 5; D = '0.05903'; E = '-0.63%'; G = '2' }
    @{ Row = 6; D = '6.642'; E = '-1.05%'; G = '2' }
    @{ Row = 7; D = '0.8540'; E = '-1.80%'; G = '2' }
    @{ Row = 8; D = '0.9417'; E = '-6.16%'; G = '2' }
    @{ Row = 9; D = '0.0006037'; E = '-94.30%'; G = '2' }
    @{ Row = 10; D = '0.1409'; E = '-0.91%'; G = '2' }
    @{ Row = 11; D = '0.04578'; E = '28.58%'; G = '2' }
    @{ Row = 12; D = '0.07089'; E = '-1.47%'; G = '2' }
    @{ Row = 13; D = '0.03136'; E = '-0.47%'; G = '2' }
    @{ Row = 14; D = '0.09158'; E = '-1.05%'; G = '2' }
    @{ Row = 15; D = '0.001534'; E = '-0.55%'; G = '2' }
    @{ Row = 16; D = '0.006221'; E = '0.79%'; G = '2' }
    @{ Row = 17; D = '3.525'; E = '0.32%'; G = '2' }
    @{ Row = 18; D = '3.193'; E = '-2.43%'; G = '2' }
    @{ Row = 19; D = '2.204'; E = '-0.09%'; G = '2' }
    @{ Row = 20; D = '0.3054'; E = $null; G = '2' }
    @{ Row = 21; D = $null; E = '-0.54%'; G = '2' }
    @{ Row = 22; D = '3.823'; E = '6.88%'; G = '2' }
    @{ Row = 23; D = '0.04264'; E = '1.53%'; G = '2' }
    @{ Row = 24; D = '0.001224'; E = '0.33%'; G = '2' }
    @{ Row = 25; D = '0.004287'; E = '-5.03%'; G = '2' }
    @{ Row = 26; D = '0.0001202'; E = '0.13%'; G = '2' }
    @{ Row = 27; D = '0.0001939'; E = '30.03%'; G = '2' }
    @{ Row = 28; D = $null; E = $null; G = '2' }
    @{ Row = 29; D = $null; E = $null; G = '2' }
    @{ Row = 30; D = $null; E = $null; G = '2' }
    @{ Row = 31; D = $null; E = $null; G = '2' }
    @{ Row = 32; D = $null; E = $null; G = '2' }
    @{ Row = 33; D = $null; E = $null; G = '2' }
    @{ Row = 34; D = $null; E = $null; G = '2' }
    @{ Row = 35; D = $null; E = $null; G = '2' }
    @{ Row = 36; D = $null; E = $null; G = '2' }
    @{ Row = 37; D = $null; E = $null; G = '2' }
    @{ Row = 38; D = $null; E = $null; G = '2' }
    @{ Row = 39; D = $null; E = $null; G = '2' }
    @{ Row = 40; D = $null; E = '-0.12%'; G = '2' }
    @{ Row = 41; D = '0.006278'; E = '57.69%'; G = '2' }
    @{ Row = 42; D = '0.1100'; E = '-0.37%'; G = '2' }
    @{ Row = 43; D = '0.002167'; E = '-7.38%'; G = '2' }
    @{ Row = 44; D = '0.01272'; E = '21.24%'; G = '2' }
    @{ Row = 45; D = '0.00005484'; E = '-0.18%'; G = '2' }
    @{ Row = 46; D = '0.00000000751'; E = '0.10%'; G = '2' }
    @{ Row = 47; D = '0.05107'; E = '-53.21%'; G = '2' }
    @{ Row = 48; D = '0.2470'; E = '11,024.02%'; G = '2' }
    @{ Row = 49; D = '0.00002103'; E = '0.10%'; G = '2' }
    @{ Row = 50; D = '0.0002003'; E = '0.10%'; G = '2' }
    @{ Row = 51; D = $null; E = $null; G = '2' }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) { $ws.Range("D$($r.Row)").Value2 = "'$($r.D)" }
    if ($null -ne $r.E) { $ws.Range("E$($r.Row)").Value2 = "'$($r.E)" }
    if ($null -ne $r.G) { $ws.Range("G$($r.Row)").Value2 = "'$($r.G)" }
}
